$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update score data for student #30 (row 32) ---
$ws.Range("D32").Value2 = 8
$ws.Range("H32").Value2 = 8
$ws.Range("I32").Value2 = 13
$ws.Range("J32").Value2 = 3

# --- Add new "Skalirano:" column (L) ---
# Header text in L2 (new shared string "Skalirano:")
# L1: standalone scaled total for the max-points row
$ws.Range("L1").Formula = "=ROUND((K1*65)/60, 1)"
# Seed L2 with the same formula pattern so a shared formula group is created
# across L2:L52, then overwrite L2 with its header text afterwards.
$ws.Range("L2").Formula = "=ROUND((K2*65)/60, 1)"
$ws.Range("L3:L52").Formula = "=ROUND((K3*65)/60, 1)"
$ws.Range("L2").Value2 = "Skalirano:"

# Column width for the new column L
$ws.Columns.Item(12).ColumnWidth = 8.38

# --- Update sheet view / selection ---
$ws.Range("M9").Select() | Out-Null
